# Update cryptos list: refresh Price (D) and Volume(1h) (E) columns for rows 2-51
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.324.09"
$ws.Range("E2").Value = "  -3.14%  "
$ws.Range("D3").Value = "1.974.03"
$ws.Range("E3").Value = "  -1.88%  "
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").Value = "'231.31"
$ws.Range("E5").Value = "  -11.34%  "
$ws.Range("D6").Value = "'0.597"
$ws.Range("E6").Value = "  -3.01%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "'53.90"
$ws.Range("E8").Value = "  -2.96%  "
$ws.Range("D9").Value = "'0.368"
$ws.Range("E9").Value = "  -3.90%  "
$ws.Range("D10").Value = "'57.88"
$ws.Range("E10").Value = "  +1.95%  "
$ws.Range("D11").Value = "'0.0748"
$ws.Range("E11").Value = "  -3.17%  "
$ws.Range("D12").Value = "'0.0981"
$ws.Range("E12").Value = "  -3.47%  "
$ws.Range("D13").Value = "2.263.42"
$ws.Range("E13").Value = "  -2.19%  "
$ws.Range("D14").Value = "'13.83"
$ws.Range("E14").Value = "  -3.25%  "
$ws.Range("D15").Value = "'19.91"
$ws.Range("E15").Value = "  -3.49%  "
$ws.Range("D16").Value = "'0.749"
$ws.Range("E16").Value = "  -6.44%  "
$ws.Range("D17").Value = "'5.01"
$ws.Range("E17").Value = "  -4.06%  "
$ws.Range("D18").Value = "1.965.40"
$ws.Range("E18").Value = "  -2.94%  "
$ws.Range("D19").Value = "36.223.60"
$ws.Range("E19").Value = "  -3.16%  "
$ws.Range("D20").Value = "'67.39"
$ws.Range("E20").Value = "  -2.88%  "
$ws.Range("D21").Value = "0.0₃0803"
$ws.Range("E21").Value = "  -4.22%  "
$ws.Range("D22").Value = "'5.20"
$ws.Range("E22").Value = "  +1.42%  "
$ws.Range("D23").Value = "'220.90"
$ws.Range("E23").Value = "  -2.88%  "
$ws.Range("E24").Value = "  +0.03%  "
$ws.Range("D25").Value = "'2.35"
$ws.Range("E25").Value = "  +1.67%  "
$ws.Range("D26").Value = "'2.33"
$ws.Range("E26").Value = "  -12.47%  "
$ws.Range("D27").Value = "'160.31"
$ws.Range("E27").Value = "  -1.74%  "
$ws.Range("D28").Value = "'8.50"
$ws.Range("E28").Value = "  -4.46%  "
$ws.Range("D29").Value = "'18.64"
$ws.Range("E29").Value = "  -4.89%  "
$ws.Range("D30").Value = "'0.124"
$ws.Range("E30").Value = "  -3.16%  "
$ws.Range("E31").Value = "  -0.43%  "
$ws.Range("D32").Value = "'0.116"
$ws.Range("E32").Value = "  -3.07%  "
$ws.Range("D33").Value = "'4.33"
$ws.Range("E33").Value = "  -5.72%  "
$ws.Range("D34").Value = "'0.0602"
$ws.Range("E34").Value = "  -7.27%  "
$ws.Range("D35").Value = "'4.23"
$ws.Range("E35").Value = "  -5.60%  "
$ws.Range("E36").Value = "  -3.04%  "
$ws.Range("E37").Value = "  -0.19%  "
$ws.Range("E38").Value = "  -3.14%  "
$ws.Range("D39").Value = "'3.22"
$ws.Range("E39").Value = "  -3.33%  "
$ws.Range("D40").Value = "'5.18"
$ws.Range("E40").Value = "  +0.02%  "
$ws.Range("D42").Value = "1.433.02"
$ws.Range("E42").Value = "  +2.63%  "
$ws.Range("D43").Value = "'0.0881"
$ws.Range("E43").Value = "  -5.91%  "
$ws.Range("D44").Value = "'0.0200"
$ws.Range("E44").Value = "  -5.93%  "
$ws.Range("D45").Value = "'1.08"
$ws.Range("E45").Value = "  -10.30%  "
$ws.Range("D46").Value = "'87.53"
$ws.Range("E46").Value = "  -1.98%  "
$ws.Range("D47").Value = "'0.985"
$ws.Range("E47").Value = "  -3.49%  "
$ws.Range("D48").Value = "'14.71"
$ws.Range("E48").Value = "  -5.69%  "
$ws.Range("D49").Value = "'2.87"
$ws.Range("E49").Value = "  -0.20%  "
$ws.Range("D50").Value = "'6.73"
$ws.Range("E50").Value = "  -4.54%  "
$ws.Range("D51").Value = "'3.57"
$ws.Range("E51").Value = "  +13.71%  "
